$d = $word.ActiveDocument

# 1. Remove "Erhebungsschwierigkeiten und daraus resultierenden "
$d.Content.Find.Execute(
    "Herausforderungen wie Erhebungsschwierigkeiten und daraus resultierenden Datenlücken",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Herausforderungen wie Datenlücken", 2) | Out-Null

# 2. "Zehnfache " -> "Zehn- "
$d.Content.Find.Execute(
    "die um das Zehnfache oder",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "die um das Zehn- oder", 2) | Out-Null

# 3. Remove " als in der Realität" before "Nach entsprechender Rücksprache", and "entsprechender " removed
$d.Content.Find.Execute(
    "niedriger ausfallen als in der Realität. Nach entsprechender Rücksprache wurden",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "niedriger ausfallen. Nach Rücksprache wurden", 2) | Out-Null

# 4. Remove " und damit die Aussagekraft auf die wahrhaftigen Daten beschränkt"
$d.Content.Find.Execute(
    "allerdings als Fehler berücksichtigt und damit die Aussagekraft auf die wahrhaftigen Daten beschränkt. Die Zusammenhänge",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "allerdings als Fehler berücksichtigt. Die Zusammenhänge", 2) | Out-Null

# 5. "unterschiedlichen " removed; reorder list to Säulen-, Linien-, Punkt- und Boxdiagramme
$d.Content.Find.Execute(
    "von unterschiedlichen Visualisierungen (z.B. Punkt-, Linien-, Säulen- und Boxdiagramme) ausgewertet",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "von Visualisierungen (z.B. Säulen-, Linien-, Punkt- und Boxdiagramme) ausgewertet", 2) | Out-Null

# 6. "mannigfaltigen " (red) -> "vielfältigen " (automatic color, not red)
$find6 = $d.Content.Find
$find6.ClearFormatting()
$find6.Text = "mannigfaltigen "
$find6.Replacement.ClearFormatting()
$find6.Replacement.Text = "vielfältigen "
$find6.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

$rng6 = $d.Content
$rng6.Find.ClearFormatting()
$rng6.Find.Text = "vielfältigen "
$rng6.Find.Execute() | Out-Null
$rng6.Font.Color = -16777216

# 8. "thematischen" -> "thematisch"
$d.Content.Find.Execute(
    "nach thematischen gewählten Indikatoren",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "nach thematisch gewählten Indikatoren", 2) | Out-Null

# 10. Add text to the first empty "Ergebnisse" paragraph (red, Arial)
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "Für die Betrachtung der Beziehungen zwischen den einzelnen Indikatoren unter ungruppierten Umständen konnten für alle jeweils tendenziell positive Zusammenhänge im schwachen bis mittelstarken Intensitätsbereich festgestellt werden. "
$p8.Range.Font.NameAscii = "Arial"
$p8.Range.Font.NameOther = "Arial"
$p8.Range.Font.NameBi = "Arial"
$p8.Range.Font.Color = 255

# 11. "Zu (i): "
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "Zu (i): "
$p9.Range.Font.NameAscii = "Arial"
$p9.Range.Font.NameOther = "Arial"
$p9.Range.Font.NameBi = "Arial"
$p9.Range.Font.Color = 255

# 12. "Zu (ii): "
$p10 = $d.Paragraphs.Item(10)
$p10.Range.Text = "Zu (ii): "
$p10.Range.Font.NameAscii = "Arial"
$p10.Range.Font.NameOther = "Arial"
$p10.Range.Font.NameBi = "Arial"
$p10.Range.Font.Color = 255

# 13. "Zu (iii):"
$p11 = $d.Paragraphs.Item(11)
$p11.Range.Text = "Zu (iii):"
$p11.Range.Font.NameAscii = "Arial"
$p11.Range.Font.NameOther = "Arial"
$p11.Range.Font.NameBi = "Arial"
$p11.Range.Font.Color = 255

# 14. Remove the last empty paragraph (merges into the placeholder paragraph),
#     then replace the placeholder text with the "Zu (iv): ..." sentence.
$p12 = $d.Paragraphs.Item(12)
$p12.Range.Delete()

$d.Content.Find.Execute(
    "[VGL. BEISPIELABGABE UND FÜGE ERGEBNISSE EIN]",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Zu (iv): der Zusammenhang ist auch im gruppierten Zustand für fast alle Ländergrößen wiederzufinden, einzig bei relativer Betrachtung haben die sehr großen Länder eine gegenläufige Entwicklung.", 2) | Out-Null

Write-Host "done part 1"
